$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels ("name" -> "Full Name", "email" -> "Email")
$ws.Range("D1").Value = "Full Name"
$ws.Range("E1").Value = "Email"

# The two hyperlink cells (E2/E3) carried two subtly different "blue link"
# fonts; unify E3's formatting onto E2's so only one such font remains.
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on D2, matching the saved cursor position.
$ws.Range("D2").Select()
